# Weekly data refresh: add the newest week's "Espárragos" (Mapocho Venta
# Directa de Santiago) observations for "Provincia de Linares" at the top
# of that block (rows 4-5), pushing all subsequent rows down by two.
#
# This mirrors what the upstream weekly ETL does: the two rows that used to
# sit at positions 4-5 (date 2022-10-13 / serial 44847) are preserved as-is,
# just shifted two rows down (to 6-7), and a brand-new pair of rows for the
# most recent week (2022-11-17 / serial 44882) is written at 4-5, re-using
# every column that did not change (market/region/category/quality/unit/
# origin/etc.) and only updating Volumen, Precio mínimo/máximo/promedio and
# Precio $/Kg.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Make room for the two new rows: shift rows 4 downward (old row 4 -> 6,
#    old row 5 -> 7, ..., old row 93 -> 95). The used range / dimension
#    updates automatically.
$ws.Rows("4:5").Insert()

# 2. Duplicate the (now shifted) old rows 4-5 -- which live at 6-7 again --
#    into the freshly inserted 4-5 so every unchanged column (A,B,C,E,F,G,H,
#    I,N,O,Q,R) is carried over exactly, including styles/number formats.
$ws.Range("A6:R7").Copy()
$ws.Range("A4").PasteSpecial()

# 3. Overwrite the cells that actually carry new data for the new week.
# New row 4 ("Primera"): fecha, volumen, precios.
$ws.Range("D4").Value = 44882
$ws.Range("J4").Value = 510
$ws.Range("K4").Value = 1200
$ws.Range("L4").Value = 1200
$ws.Range("M4").Value = 1200
$ws.Range("P4").Value = 1200

# New row 5 ("Segunda"): fecha, volumen, precios.
$ws.Range("D5").Value = 44882
$ws.Range("J5").Value = 450
$ws.Range("K5").Value = 800
$ws.Range("L5").Value = 800
$ws.Range("M5").Value = 800
$ws.Range("P5").Value = 800
